# Applies the cryptos.xlsx "Updated symbol list" edit: refreshed Price (D)
# and Volume(1h) (E) quotes for the rows whose figures moved in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col=4; Text="291.00"},
    @{Row=2; Col=5; Text="-0.52%"},
    @{Row=3; Col=4; Text="30.75"},
    @{Row=3; Col=5; Text="-2.01%"},
    @{Row=4; Col=4; Text="4.891"},
    @{Row=4; Col=5; Text="-1.51%"},
    @{Row=5; Col=4; Text="0.07262"},
    @{Row=5; Col=5; Text="-1.55%"},
    @{Row=6; Col=4; Text="2.353"},
    @{Row=6; Col=5; Text="29.68%"},
    @{Row=7; Col=4; Text="7.662"},
    @{Row=7; Col=5; Text="0.00%"},
    @{Row=8; Col=5; Text="-1.32%"},
    @{Row=9; Col=4; Text="0.8973"},
    @{Row=9; Col=5; Text="-1.58%"},
    @{Row=10; Col=4; Text="0.1670"},
    @{Row=10; Col=5; Text="1.41%"},
    @{Row=11; Col=4; Text="0.07974"},
    @{Row=11; Col=5; Text="4.91%"},
    @{Row=12; Col=4; Text="0.08188"},
    @{Row=12; Col=5; Text="-0.02%"},
    @{Row=13; Col=4; Text="0.03082"},
    @{Row=13; Col=5; Text="3.10%"},
    @{Row=14; Col=4; Text="0.1003"},
    @{Row=14; Col=5; Text="0.71%"},
    @{Row=15; Col=4; Text="0.001497"},
    @{Row=15; Col=5; Text="-0.30%"},
    @{Row=16; Col=4; Text="0.005828"},
    @{Row=16; Col=5; Text="3.32%"},
    @{Row=17; Col=4; Text="3.474"},
    @{Row=17; Col=5; Text="-0.01%"},
    @{Row=18; Col=4; Text="2.078"},
    @{Row=18; Col=5; Text="-2.31%"},
    @{Row=19; Col=5; Text="0.79%"},
    @{Row=20; Col=4; Text="0.1298"},
    @{Row=20; Col=5; Text="0.35%"},
    @{Row=21; Col=4; Text="3.970"},
    @{Row=21; Col=5; Text="-8.15%"},
    @{Row=22; Col=4; Text="0.2300"},
    @{Row=22; Col=5; Text="16.39%"},
    @{Row=23; Col=4; Text="0.04523"},
    @{Row=23; Col=5; Text="0.54%"},
    @{Row=24; Col=5; Text="-1.15%"},
    @{Row=25; Col=4; Text="0.004414"},
    @{Row=25; Col=5; Text="8.95%"},
    @{Row=26; Col=4; Text="0.0001301"},
    @{Row=26; Col=5; Text="4.05%"},
    @{Row=27; Col=4; Text="0.0003392"},
    @{Row=27; Col=5; Text="-95.48%"},
    @{Row=39; Col=4; Text="0.01589"},
    @{Row=39; Col=5; Text="-2.96%"},
    @{Row=40; Col=4; Text="0.04378"},
    @{Row=40; Col=5; Text="-0.33%"},
    @{Row=41; Col=4; Text="0.007321"},
    @{Row=41; Col=5; Text="-1.57%"},
    @{Row=43; Col=4; Text="0.1316"},
    @{Row=44; Col=4; Text="0.002017"},
    @{Row=44; Col=5; Text="-2.25%"},
    @{Row=45; Col=4; Text="0.009505"},
    @{Row=45; Col=5; Text="-13.96%"},
    @{Row=46; Col=5; Text="-4.82%"},
    @{Row=47; Col=4; Text="0.00000000751"},
    @{Row=47; Col=5; Text="0.02%"},
    @{Row=48; Col=5; Text="17.87%"},
    @{Row=49; Col=4; Text="0.002896"},
    @{Row=49; Col=5; Text="-3.50%"},
    @{Row=50; Col=4; Text="0.00002102"},
    @{Row=50; Col=5; Text="0.02%"},
    @{Row=51; Col=4; Text="0.0002002"},
    @{Row=51; Col=5; Text="0.02%"}
)

foreach ($u in $updates) {
    # Force text formatting first so numeric-looking quotes (e.g. "291.00")
    # and percentages (e.g. "-0.52%") are kept as literal strings, matching
    # how this sheet stores every Price/Volume cell.
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Text
}
